$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Weights RP")

# Fix weights for existing representative periods (rp01-rp04): 48 -> 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1

# Add two new representative periods rp05 and rp06 with weight 1
$ws.Range("B11").Value = "rp05"
$ws.Range("C11").Value = 1

$ws.Range("B12").Value = "rp06"
$ws.Range("C12").Value = 1

# Copy style from existing rows (B10/C10) down to new rows (B11:C12)
$ws.Range("B10:C10").Copy() | Out-Null
$ws.Range("B11:C12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update selection to match final state
$ws.Range("D15").Select() | Out-Null
